$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the "Pin Size" value for row 6 (B6)
$ws.Range("B6").Value = 1.57

# Overwrite the previously-shared formula in C6 with a hard-coded value
$ws.Range("C6").Value = 1.6

# Move the active selection to C7 (as if the user had just tabbed/entered down)
$ws.Activate()
$ws.Range("C7").Select()
